# feat: add 2022-Q1 data
#
# Plan (matches the OOXML diff):
#   1. The existing "总计" (summary) sheet (2nd sheet) is renamed to "2022-Q1"
#      and its content is replaced by the 2022-Q1 per-fund holdings table
#      (same shape/style as the existing "2021-Q4" sheet).
#   2. A brand-new "总计" sheet is appended at the end, with the summary
#      table now covering both quarters (2022-Q1 on top, 2021-Q4 below).

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item(1)      # "2021-Q4" - untouched, used as a style template
$wsQ1 = $wb.Worksheets.Item(2)      # currently "总计" - becomes "2022-Q1"

# ---------------------------------------------------------------------
# 1) Turn the old "总计" sheet into the "2022-Q1" per-fund holdings sheet
# ---------------------------------------------------------------------

# Clear whatever was there before (old summary table: A1:D2)
$wsQ1.Cells.Clear()

$wsQ1.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

$fundRows = @(
    @(0, "501305", "汇添富中证港股通高股息投资指数（LOF）A", "1.59", "93.08", "4.61", "0.0733", 4),
    @(1, "004532", "民生加银中证港股通高股息精选指数A",       "0.26", "94.88", "9.04", "0.0235", 1),
    @(2, "501306", "汇添富中证港股通高股息投资指数（LOF）C", "0.21", "93.08", "4.61", "0.0097", 4),
    @(3, "004533", "民生加银中证港股通高股息精选指数C",       "0.10", "94.88", "9.04", "0.0090", 1),
    @(4, "501307", "银河中证沪港深高股息指数（LOF）A",         "0.19", "91.35", "1.81", "0.0034", 4),
    @(5, "501308", "银河中证沪港深高股息指数（LOF）C",         "0.01", "91.35", "1.81", "0.0002", 4)
)

# Header row (B1:H1) - force text so numeric-looking headers aren't coerced,
# though none of these look numeric; kept uniform with the data-cell approach.
$headerRange = $wsQ1.Range("B1:H1")
$headerRange.NumberFormat = "@"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ1.Cells.Item(1, 2 + $i).Value = $headers[$i]
}
$headerRange.ClearFormats()

# Data rows starting at row 2. Columns B..G are text-like (kept as typed
# strings, even the numeric-looking ones), column A and H are real numbers.
for ($r = 0; $r -lt $fundRows.Length; $r++) {
    $row = 2 + $r
    $values = $fundRows[$r]

    $wsQ1.Cells.Item($row, 1).Value = $values[0]          # A: index number

    $textRange = $wsQ1.Range("B$row`:G$row")
    $textRange.NumberFormat = "@"
    $wsQ1.Cells.Item($row, 2).Value = $values[1]          # B: 基金代码
    $wsQ1.Cells.Item($row, 3).Value = $values[2]          # C: 基金名称
    $wsQ1.Cells.Item($row, 4).Value = $values[3]          # D: 基金规模
    $wsQ1.Cells.Item($row, 5).Value = $values[4]          # E: 股票总仓位
    $wsQ1.Cells.Item($row, 6).Value = $values[5]          # F: 仓位占比
    $wsQ1.Cells.Item($row, 7).Value = $values[6]          # G: 持有市值(亿元)
    $textRange.ClearFormats()

    $wsQ1.Cells.Item($row, 8).Value = $values[7]          # H: 仓位排名 (number)
}

# Re-apply the same look as the "2021-Q4" sheet: bold+bordered header row
# and bold+bordered index column (A).
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)

$wsQ4.Range("A2:A7").Copy()
$wsQ1.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Append a new "总计" sheet with the combined summary table
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTotal = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsTotal.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
$totalHeaderRange = $wsTotal.Range("B1:D1")
$totalHeaderRange.NumberFormat = "@"
for ($i = 0; $i -lt $totalHeaders.Length; $i++) {
    $wsTotal.Cells.Item(1, 2 + $i).Value = $totalHeaders[$i]
}
$totalHeaderRange.ClearFormats()

# Row 2: 2022-Q1
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Range("B2").NumberFormat = "@"
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Range("B2").ClearFormats()
$wsTotal.Cells.Item(2, 3).Value = 6
$wsTotal.Cells.Item(2, 4).Value = 0.12

# Row 3: 2021-Q4
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Range("B3").NumberFormat = "@"
$wsTotal.Cells.Item(3, 2).Value = "2021-Q4"
$wsTotal.Range("B3").ClearFormats()
$wsTotal.Cells.Item(3, 3).Value = 6
$wsTotal.Cells.Item(3, 4).Value = 0.66

# Apply the same header / index-column look as the other sheets.
$wsQ4.Range("B1:D1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)

$wsQ4.Range("A2:A3").Copy()
$wsTotal.Range("A2:A3").PasteSpecial(-4122)
